$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before row 81 - this shifts the existing rows 81..94 down to 82..95
$ws.Rows(81).Insert()

# Populate the newly inserted row 81 with the new weekly record
$ws.Range("A81").Value = 1
$ws.Range("B81").Value = "Agrícola del Norte S.A. de Arica"
$ws.Range("C81").Value = "Arica y Parinacota"
$ws.Range("D81").Value = 44776
$ws.Range("D81").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("E81").Value = 15
$ws.Range("F81").Value = "Fruta"
$ws.Range("G81").Value = 100102
$ws.Range("H81").Value = "Cítricos"
$ws.Range("I81").Value = 100102005
$ws.Range("J81").Value = "Naranja"
$ws.Range("K81").Value = "Fukumoto"
$ws.Range("L81").Value = "Segunda"
$ws.Range("M81").Value = 270
$ws.Range("N81").Value = 500
$ws.Range("O81").Value = 550
$ws.Range("P81").Value = 525
$ws.Range("Q81").Value = "`$/kilo (en caja de 20 kilos)"
$ws.Range("R81").Value = "Región de Coquimbo"
$ws.Range("S81").Value = 525
$ws.Range("T81").Value = 1
